$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.736.07"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "1.636.38"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'213.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "'0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "'0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'19.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "'0.0834"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").Value = "1.860.38"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.607.60"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "26.688.17"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "'63.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'209.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'4.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'9.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "'6.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'145.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'15.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").Value = "'6.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "'0.0520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.77%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "1.165.98"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "'0.0167"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "'0.816"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").Value = "'0.787"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "1.770.44"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").Value = "'92.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").Value = "'1.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "'54.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("D49").Value = "'7.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.71%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0512"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.86%  "
